$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 191.66667
$ws.Range("I12").Value = 137.5
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 137.5
$ws.Range("L12").Value = 300
$ws.Range("M12").Value = 32.5
$ws.Range("N12").Value = -640

$ws.Range("H21").Value = 30399.9
$ws.Range("J21").Value = 31555.555
$ws.Range("L21").Value = 31555.555
$ws.Range("N21").Value = -32491.555

$ws.Range("H23").Value = 30399.9
$ws.Range("J23").Value = 31555.555
$ws.Range("L23").Value = 31555.555
$ws.Range("N23").Value = -32023.555

$ws.Range("H132").Value = 1509.1892
$ws.Range("I132").Value = 994.8387
$ws.Range("J132").Value = 4166.6665
$ws.Range("K132").Value = 2984.5161
$ws.Range("L132").Value = 12499.9995
$ws.Range("M132").Value = -454.5160999999998
$ws.Range("N132").Value = -17559.9995

$ws.Range("H135").Value = 3750.6667
$ws.Range("I135").Value = 4666.609
$ws.Range("K135").Value = 41999.481
$ws.Range("M135").Value = -39464.481

$ws.Range("H138").Value = 4135.4414
$ws.Range("I138").Value = 3135.087
$ws.Range("J138").Value = 6227.091
$ws.Range("K138").Value = 9405.261
$ws.Range("L138").Value = 18681.273
$ws.Range("M138").Value = -4265.261
$ws.Range("N138").Value = -28961.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6599.9873
$ws.Range("I32").Value = 2834.6719
$ws.Range("J32").Value = 21661.25
$ws.Range("K32").Value = 2834.6719
$ws.Range("L32").Value = 21661.25
$ws.Range("M32").Value = -2547.6719
$ws.Range("N32").Value = -22235.25

$ws.Range("H107").Value = 60266.668
$ws.Range("J107").Value = 60266.668
$ws.Range("L107").Value = 60266.668
$ws.Range("N107").Value = -67946.66800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1696.7693
$ws.Range("I86").Value = 1568.591
$ws.Range("J86").Value = 2401.75
$ws.Range("K86").Value = 1568.591
$ws.Range("L86").Value = 2401.75
$ws.Range("M86").Value = -445.5909999999999
$ws.Range("N86").Value = -4647.75

$ws.Range("H89").Value = 1696.7693
$ws.Range("I89").Value = 1568.591
$ws.Range("J89").Value = 2401.75
$ws.Range("K89").Value = 7842.955
$ws.Range("L89").Value = 12008.75
$ws.Range("M89").Value = -2226.955
$ws.Range("N89").Value = -23240.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3762.5
$ws.Range("I132").Value = 3480.5
$ws.Range("K132").Value = 10441.5
$ws.Range("M132").Value = -7911.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 1870
$ws.Range("I31").Value = 382.22223
$ws.Range("J31").Value = 6333.3335
$ws.Range("K31").Value = 1146.66669
$ws.Range("L31").Value = 19000.0005
$ws.Range("M31").Value = -858.66669
$ws.Range("N31").Value = -19576.0005

$ws.Range("H49").Value = 3825.5
$ws.Range("I49").Value = 1953
$ws.Range("J49").Value = 4200
$ws.Range("K49").Value = 5859
$ws.Range("L49").Value = 12600
$ws.Range("M49").Value = -5703
$ws.Range("N49").Value = -12912

$ws.Range("H54").Value = 2611.5
$ws.Range("I54").Value = 2011.8182
$ws.Range("J54").Value = 3344.4443
$ws.Range("K54").Value = 6035.4546
$ws.Range("L54").Value = 10033.3329
$ws.Range("M54").Value = -5476.4546
$ws.Range("N54").Value = -11151.3329

$ws.Range("H62").Value = 2941.6667
$ws.Range("I62").Value = 475
$ws.Range("J62").Value = 4175
$ws.Range("K62").Value = 1425
$ws.Range("L62").Value = 12525
$ws.Range("M62").Value = -739
$ws.Range("N62").Value = -13897

$ws.Range("H63").Value = 5077.6
$ws.Range("I63").Value = 3150
$ws.Range("J63").Value = 5559.5
$ws.Range("K63").Value = 9450
$ws.Range("L63").Value = 16678.5
$ws.Range("M63").Value = -8701
$ws.Range("N63").Value = -18176.5

$ws.Range("H65").Value = 2941.6667
$ws.Range("I65").Value = 475
$ws.Range("J65").Value = 4175
$ws.Range("K65").Value = 4275
$ws.Range("L65").Value = 37575
$ws.Range("M65").Value = -843
$ws.Range("N65").Value = -44439

$ws.Range("H66").Value = 5077.6
$ws.Range("I66").Value = 3150
$ws.Range("J66").Value = 5559.5
$ws.Range("K66").Value = 28350
$ws.Range("L66").Value = 50035.5
$ws.Range("M66").Value = -24606
$ws.Range("N66").Value = -57523.5

$ws.Range("H68").Value = 1358.0116
$ws.Range("I68").Value = 958.9545000000001
$ws.Range("J68").Value = 1776.0714
$ws.Range("K68").Value = 2876.8635
$ws.Range("L68").Value = 5328.2142
$ws.Range("M68").Value = -2065.8635
$ws.Range("N68").Value = -6950.2142

$ws.Range("H71").Value = 1358.0116
$ws.Range("I71").Value = 958.9545000000001
$ws.Range("J71").Value = 1776.0714
$ws.Range("K71").Value = 8630.5905
$ws.Range("L71").Value = 15984.6426
$ws.Range("M71").Value = -4574.5905
$ws.Range("N71").Value = -24096.6426

$ws.Range("H74").Value = 7079.7
$ws.Range("I74").Value = 600
$ws.Range("J74").Value = 9856.714
$ws.Range("K74").Value = 1800
$ws.Range("L74").Value = 29570.142
$ws.Range("M74").Value = -739
$ws.Range("N74").Value = -31692.142

$ws.Range("H77").Value = 7079.7
$ws.Range("I77").Value = 600
$ws.Range("J77").Value = 9856.714
$ws.Range("K77").Value = 5400
$ws.Range("L77").Value = 88710.42600000001
$ws.Range("M77").Value = -96
$ws.Range("N77").Value = -99318.42600000001

$ws.Range("H92").Value = 476.66666
$ws.Range("I92").Value = 528.5714
$ws.Range("J92").Value = 295
$ws.Range("K92").Value = 1585.7142
$ws.Range("L92").Value = 885
$ws.Range("M92").Value = -337.7142000000001
$ws.Range("N92").Value = -3381

$ws.Range("H93").Value = 7866.6665
$ws.Range("J93").Value = 7866.6665
$ws.Range("L93").Value = 23599.9995
$ws.Range("N93").Value = -27343.9995

$ws.Range("H94").Value = 6601.4546
$ws.Range("J94").Value = 7624
$ws.Range("L94").Value = 22872
$ws.Range("N94").Value = -24224

$ws.Range("H96").Value = 5799.9
$ws.Range("J96").Value = 5799.9
$ws.Range("L96").Value = 17399.7
$ws.Range("N96").Value = -21517.7

$ws.Range("H97").Value = 496.42856
$ws.Range("I97").Value = 275
$ws.Range("J97").Value = 533.3333
$ws.Range("K97").Value = 825
$ws.Range("L97").Value = 1599.9999
$ws.Range("M97").Value = -329
$ws.Range("N97").Value = -2591.9999

$ws.Range("H98").Value = 1331.7778
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 1331.7778
$ws.Range("K98").Value = 0
$ws.Range("L98").ClearContents()
$ws.Range("M98").Value = 3995.3334
$ws.Range("N98").Value = -6991.3334

$ws.Range("H99").Value = 2064.6667
$ws.Range("I99").Value = 1422
$ws.Range("J99").Value = 3350
$ws.Range("K99").Value = 4266
$ws.Range("L99").Value = 10050
$ws.Range("M99").Value = -2020
$ws.Range("N99").Value = -14542

$ws.Range("H101").Value = 8905.799999999999
$ws.Range("J101").Value = 8905.799999999999
$ws.Range("L101").Value = 26717.4
$ws.Range("N101").Value = -31585.4

$ws.Range("H102").Value = 3500
$ws.Range("J102").Value = 3500
$ws.Range("L102").Value = 10500
$ws.Range("N102").Value = -15368

$ws.Range("H106").Value = 8905.883
$ws.Range("J106").Value = 8905.883
$ws.Range("L106").Value = 26717.649
$ws.Range("N106").Value = -28609.649

$ws.Range("H108").Value = 354.5
$ws.Range("I108").Value = 354.5
$ws.Range("K108").Value = 1063.5
$ws.Range("M108").Value = 1816.5

$ws.Range("H114").Value = 1389.2
$ws.Range("I114").Value = 695.4
$ws.Range("J114").Value = 2083
$ws.Range("K114").Value = 2086.2
$ws.Range("L114").Value = 6249
$ws.Range("M114").Value = 1167.8
$ws.Range("N114").Value = -12757

$ws.Range("H116").Value = 604.2857
$ws.Range("I116").Value = 604.2857
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1812.8571
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = 1629.1429

$ws.Range("H118").Value = 2528070.5
$ws.Range("I118").Value = 2120
$ws.Range("J118").Value = 4633029.5
$ws.Range("K118").Value = 6360
$ws.Range("L118").Value = 13899088.5
$ws.Range("M118").Value = -5117
$ws.Range("N118").Value = -13901574.5

$ws.Range("H130").Value = 4401.8
$ws.Range("I130").Value = 2015
$ws.Range("J130").Value = 5993
$ws.Range("K130").Value = 6045
$ws.Range("L130").Value = 17979
$ws.Range("M130").Value = -1025
$ws.Range("N130").Value = -28019

$ws.Range("H137").Value = 1436704.2
$ws.Range("I137").Value = 8758.462
$ws.Range("J137").Value = 20000000
$ws.Range("K137").Value = 26275.386
$ws.Range("L137").Value = 60000000
$ws.Range("M137").Value = -21175.386
$ws.Range("N137").Value = -60010200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5266.364
$ws.Range("I132").Value = 7330
$ws.Range("J132").Value = 2465.7144
$ws.Range("K132").Value = 21990
$ws.Range("L132").Value = 7397.1432
$ws.Range("M132").Value = -19460
$ws.Range("N132").Value = -12457.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 1251.5
$ws.Range("I7").Value = 1251.5
$ws.Range("K7").Value = 1251.5
$ws.Range("M7").Value = -1138.5

$ws.Range("H132").Value = 1122.069
$ws.Range("I132").Value = 745.26086
$ws.Range("J132").Value = 2566.5
$ws.Range("K132").Value = 2235.78258
$ws.Range("L132").Value = 7699.5
$ws.Range("M132").Value = 294.2174199999999
$ws.Range("N132").Value = -12759.5
